$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. Using NumberFormat "@" (Text) while
# assigning forces Excel to keep numeric-looking strings (e.g. "202.85", "1.00")
# as literal text instead of coercing them to numbers, matching the source data
# which stores these as plain strings. ClearFormats() afterwards removes the
# temporary text-format style so the cell keeps its original (default) styling.
$updates = [ordered]@{
    'D2' = '76.489.42'
    'E2' = '  +0.26%  '
    'D3' = '3.050.83'
    'E3' = '  +4.23%  '
    'E4' = '  +0.04%  '
    'D5' = '202.85'
    'E5' = '  -0.45%  '
    'D6' = '625.61'
    'E6' = '  +4.43%  '
    'D7' = '1.00'
    'E7' = '  +0.06%  '
    'E8' = '  +0.06%  '
    'D9' = '0.208'
    'E9' = '  +5.50%  '
    'D10' = '3.049.82'
    'E10' = '  +4.18%  '
    'D11' = '0.439'
    'E11' = '  +1.48%  '
    'E12' = '  -0.61%  '
    'E13' = '  +4.76%  '
    'D14' = '3.612.74'
    'E14' = '  +4.28%  '
    'D15' = '29.53'
    'E15' = '  +5.32%  '
    'D16' = '76.455.80'
    'E16' = '  +0.40%  '
    'E17' = '  +1.45%  '
    'D18' = '3.045.91'
    'E18' = '  +4.26%  '
    'D19' = '13.50'
    'E19' = '  +4.09%  '
    'D20' = '9.10'
    'E20' = '  +2.82%  '
    'D21' = '374.73'
    'E21' = '  -0.09%  '
    'E22' = '  -0.50%  '
    'E23' = '  +0.97%  '
    'D24' = '73.59'
    'E24' = '  +2.84%  '
    'D25' = '3.208.49'
    'E25' = '  +5.06%  '
    'E26' = '  +4.04%  '
    'E27' = '  -0.11%  '
    'D28' = '9.95'
    'E28' = '  +2.25%  '
    'E29' = '  +2.38%  '
    'D30' = '0.996'
    'E30' = '  -0.50%  '
    'D31' = '8.31'
    'E31' = '  +7.18%  '
    'E32' = '  +0.64%  '
    'D33' = '509.98'
    'E33' = '  +0.72%  '
    'E34' = '  +6.66%  '
    'D35' = '1.00'
    'E35' = '  -0.02%  '
    'D36' = '20.93'
    'E36' = '  +3.14%  '
    'D37' = '163.09'
    'E37' = '  -0.50%  '
    'D38' = '0.388'
    'E38' = '  +8.14%  '
    'E39' = '  +2.05%  '
    'B40' = 'Aave'
    'C40' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D40' = '190.59'
    'E40' = '  +4.68%  '
    'B41' = 'Cronos'
    'C41' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D41' = '0.106'
    'E41' = '  +0.15%  '
    'B42' = 'Mantle'
    'C42' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    'D42' = '0.850'
    'E42' = '  +29.43%  '
    'E43' = '  +0.21%  '
    'E44' = '  +0.00%  '
    'E45' = '  +3.83%  '
    'D46' = '1.33'
    'E46' = '  +11.03%  '
    'D47' = '42.14'
    'E47' = '  +5.39%  '
    'E48' = '  -0.20%  '
    'D49' = '2.45'
    'E49' = '  +3.49%  '
    'D50' = '0.611'
    'E50' = '  +6.50%  '
    'E51' = '  +4.63%  '
}

foreach ($addr in $updates.Keys) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $updates[$addr]
    $c.ClearFormats()
}
